$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the flight schedule log (rows 947-963).
# Columns: A = date serial, B = scheduled flights, C = actual flights,
# D = C/B (percentage completed), carried down as a formula like the
# existing rows above.
$data = @(
    @(947, 44875, 81, 79),
    @(948, 44876, 83, 78),
    @(949, 44877, 43, 43),
    @(950, 44878, 67, 65),
    @(951, 44879, 63, 61),
    @(952, 44880, 60, 59),
    @(953, 44881, 66, 64),
    @(954, 44882, 87, 84),
    @(955, 44883, 78, 75),
    @(956, 44884, 52, 52),
    @(957, 44885, 65, 64),
    @(958, 44886, 58, 58),
    @(959, 44887, 63, 62),
    @(960, 44888, 62, 62),
    @(961, 44889, 39, 36),
    @(962, 44890, 44, 43),
    @(963, 44891, 62, 61)
)

$firstRow = $data[0][0]
$lastRow = $data[$data.Count - 1][0]

# Copy the formatting (number formats / styles) of the last existing row
# down across the new rows before filling in the actual values.
$ws.Range("A946:D946").Copy()
$ws.Range("A" + $firstRow + ":D" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Formula = "=C$r/B$r"
}

# Restore the view state to match the saved selection/scroll position.
$ws.Range("G959").Select()
$excel.ActiveWindow.ScrollRow = 929
$excel.ActiveWindow.ScrollColumn = 1
